$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.267.49'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '1.618.55'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '211.98'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("E10").Value = '  +5.20%  '
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '1.843.01'
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").Value = '1.620.78'
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").Value = '26.272.67'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '62.26'
$ws.Range("E17").Value = '  +3.90%  '
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '201.04'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  +2.95%  '
$ws.Range("D25").Value = '144.09'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("E30").Value = '  +9.04%  '
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("E32").Value = '  +1.85%  '
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("D36").Value = '1.177.13'
$ws.Range("E36").Value = '  +4.86%  '
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("D38").Value = '0.801'
$ws.Range("E38").Value = '  +2.97%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("E43").Value = '  +5.02%  '
$ws.Range("D44").Value = '1.754.79'
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("D45").Value = '92.59'
$ws.Range("E45").Value = '  +1.03%  '
$ws.Range("E46").Value = '  +13.86%  '
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").Value = '53.62'
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("E51").Value = '  -0.19%  '
